$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (above the current first data row),
# shifting all existing data rows down by one.
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above (the header);
# strip that so the new data row matches the plain (unstyled) data rows.
$ws.Range("A2:Q2").ClearFormats()

# Populate the new row's values. Force text number format on the
# date-like columns first so Excel stores them as literal text (shared
# strings) instead of auto-converting them to date serials, then strip
# the format again so no stray style survives on the cell.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2024-04-24"

$ws.Cells.Item(2, 2).Value = "하나33호스팩"
$ws.Cells.Item(2, 3).Value = "코스닥"
$ws.Cells.Item(2, 4).Value = 70
$ws.Cells.Item(2, 5).Value = "하나"
$ws.Cells.Item(2, 6).Value = 70
$ws.Cells.Item(2, 7).Value = "-"
$ws.Cells.Item(2, 8).Value = "-"
$ws.Cells.Item(2, 9).Value = "-"
$ws.Cells.Item(2, 10).Value = "-"
$ws.Cells.Item(2, 11).Value = "대표"
$ws.Cells.Item(2, 12).Value = "-"
$ws.Cells.Item(2, 13).Value = 2000
$ws.Cells.Item(2, 14).Value = 100

$ws.Cells.Item(2, 15).NumberFormat = "@"
$ws.Cells.Item(2, 15).Value = "2024-04-15"

$ws.Cells.Item(2, 16).NumberFormat = "@"
$ws.Cells.Item(2, 16).Value = "2024-04-18"

$ws.Cells.Item(2, 17).Value = 2625000

# Drop any styling picked up while typing the text-formatted date cells
# so the new row ends up with no style overrides, matching the other
# plain data rows.
$ws.Range("A2:Q2").ClearFormats()
